$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 907, shifting the existing
# rows 907:946 down to 911:950 (dimension grows from R946 to R950).
$ws.Rows("907:910").Insert()

# --- New row 907 ---
$ws.Range("A907").Value = 6
$ws.Range("B907").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C907").Value = "Metropolitana"
$ws.Range("D907").Value = 44509
$ws.Range("E907").Value = 13
$ws.Range("F907").Value = 100114013
$ws.Range("G907").Value = "Zanahoria"
$ws.Range("H907").Value = "Sin especificar"
$ws.Range("I907").Value = "Primera"
$ws.Range("J907").Value = 2040
$ws.Range("K907").Value = 8000
$ws.Range("L907").Value = 9000
$ws.Range("M907").Value = 8328
$ws.Range("N907").Value = "$/saco 20 kilos"
$ws.Range("O907").Value = "Chillán"
$ws.Range("P907").Value = 416
$ws.Range("Q907").Value = 20
$ws.Range("R907").Value = "Hortaliza"

# --- New row 908 ---
$ws.Range("A908").Value = 6
$ws.Range("B908").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C908").Value = "Metropolitana"
$ws.Range("D908").Value = 44509
$ws.Range("E908").Value = 13
$ws.Range("F908").Value = 100114013
$ws.Range("G908").Value = "Zanahoria"
$ws.Range("H908").Value = "Sin especificar"
$ws.Range("I908").Value = "Primera"
$ws.Range("J908").Value = 1650
$ws.Range("K908").Value = 7500
$ws.Range("L908").Value = 8000
$ws.Range("M908").Value = 7794
$ws.Range("N908").Value = "$/saco 20 kilos"
$ws.Range("O908").Value = "Región Metropolitana"
$ws.Range("P908").Value = 390
$ws.Range("Q908").Value = 20
$ws.Range("R908").Value = "Hortaliza"

# --- New row 909 ---
$ws.Range("A909").Value = 6
$ws.Range("B909").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C909").Value = "Metropolitana"
$ws.Range("D909").Value = 44509
$ws.Range("E909").Value = 13
$ws.Range("F909").Value = 100114013
$ws.Range("G909").Value = "Zanahoria"
$ws.Range("H909").Value = "Sin especificar"
$ws.Range("I909").Value = "Segunda"
$ws.Range("J909").Value = 470
$ws.Range("K909").Value = 6000
$ws.Range("L909").Value = 6000
$ws.Range("M909").Value = 6000
$ws.Range("N909").Value = "$/saco 20 kilos"
$ws.Range("O909").Value = "Chillán"
$ws.Range("P909").Value = 300
$ws.Range("Q909").Value = 20
$ws.Range("R909").Value = "Hortaliza"

# --- New row 910 ---
$ws.Range("A910").Value = 6
$ws.Range("B910").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C910").Value = "Metropolitana"
$ws.Range("D910").Value = 44509
$ws.Range("E910").Value = 13
$ws.Range("F910").Value = 100114013
$ws.Range("G910").Value = "Zanahoria"
$ws.Range("H910").Value = "Sin especificar"
$ws.Range("I910").Value = "Segunda"
$ws.Range("J910").Value = 350
$ws.Range("K910").Value = 6000
$ws.Range("L910").Value = 6000
$ws.Range("M910").Value = 6000
$ws.Range("N910").Value = "$/saco 20 kilos"
$ws.Range("O910").Value = "Región Metropolitana"
$ws.Range("P910").Value = 300
$ws.Range("Q910").Value = 20
$ws.Range("R910").Value = "Hortaliza"
